# Add new columns I (I0) and J (IF) to the Hader Josh save-data sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cells (row 1) ---
# Copy formatting (style) from the existing header cell H1 onto the new
# header cells I1:J1, then set their text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# --- Data rows (rows 2-65) ---
$data = @(
    @{Row=2; I=9; J=9},
    @{Row=3; I=9; J=9},
    @{Row=4; I=9; J=9},
    @{Row=5; I=8; J=9},
    @{Row=6; I=9; J=9},
    @{Row=7; I=9; J=9},
    @{Row=8; I=9; J=9},
    @{Row=9; I=9; J=9},
    @{Row=10; I=9; J=9},
    @{Row=11; I=9; J=9},
    @{Row=12; I=9; J=9},
    @{Row=13; I=9; J=9},
    @{Row=14; I=9; J=9},
    @{Row=15; I=9; J=9},
    @{Row=16; I=9; J=9},
    @{Row=17; I=9; J=9},
    @{Row=18; I=9; J=9},
    @{Row=19; I=8; J=8},
    @{Row=20; I=8; J=8},
    @{Row=21; I=9; J=9},
    @{Row=22; I=9; J=9},
    @{Row=23; I=9; J=9},
    @{Row=24; I=9; J=9},
    @{Row=25; I=9; J=9},
    @{Row=26; I=9; J=9},
    @{Row=27; I=9; J=9},
    @{Row=28; I=9; J=9},
    @{Row=29; I=9; J=9},
    @{Row=30; I=9; J=9},
    @{Row=31; I=9; J=9},
    @{Row=32; I=9; J=9},
    @{Row=33; I=9; J=9},
    @{Row=34; I=9; J=9},
    @{Row=35; I=9; J=9},
    @{Row=36; I=9; J=9},
    @{Row=37; I=9; J=9},
    @{Row=38; I=9; J=9},
    @{Row=39; I=9; J=9},
    @{Row=40; I=9; J=9},
    @{Row=41; I=9; J=9},
    @{Row=42; I=9; J=9},
    @{Row=43; I=9; J=9},
    @{Row=44; I=9; J=9},
    @{Row=45; I=9; J=9},
    @{Row=46; I=9; J=9},
    @{Row=47; I=9; J=9},
    @{Row=48; I=9; J=9},
    @{Row=49; I=9; J=9},
    @{Row=50; I=9; J=9},
    @{Row=51; I=9; J=9},
    @{Row=52; I=9; J=9},
    @{Row=53; I=9; J=9},
    @{Row=54; I=9; J=9},
    @{Row=55; I=9; J=9},
    @{Row=56; I=9; J=9},
    @{Row=57; I=9; J=9},
    @{Row=58; I=9; J=9},
    @{Row=59; I=9; J=9},
    @{Row=60; I=9; J=9},
    @{Row=61; I=9; J=9},
    @{Row=62; I=7; J=7},
    @{Row=63; I=6; J=6},
    @{Row=64; I=6; J=6},
    @{Row=65; I=4; J=4}
)

foreach ($entry in $data) {
    $r = $entry.Row
    $ws.Cells.Item($r, 9).Value = $entry.I   # column I
    $ws.Cells.Item($r, 10).Value = $entry.J  # column J
}
